$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2360.7083
$ws.Range("J121").Value = 2423.1365
$ws.Range("L121").Value = 7269.4095
$ws.Range("N121").Value = -10763.4095

$ws.Range("H132").Value = 14928327
$ws.Range("I132").Value = 15387538
$ws.Range("K132").Value = 46162614
$ws.Range("M132").Value = -46160084

$ws.Range("H137").Value = 99963.61
$ws.Range("I137").Value = 198461.56
$ws.Range("J137").Value = 1465.6666
$ws.Range("K137").Value = 595384.6799999999
$ws.Range("L137").Value = 4396.9998
$ws.Range("M137").Value = -592834.6799999999
$ws.Range("N137").Value = -9496.9998

$ws.Range("H138").Value = 4921.448
$ws.Range("J138").Value = 4944.9805
$ws.Range("L138").Value = 14834.9415
$ws.Range("N138").Value = -25114.9415

$ws.Range("H141").Value = 1621.875
$ws.Range("I141").Value = 1658.5714
$ws.Range("J141").Value = 1365
$ws.Range("K141").Value = 4975.7142
$ws.Range("L141").Value = 4095
$ws.Range("M141").Value = 204.2857999999997
$ws.Range("N141").Value = -14455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4262.83
$ws.Range("I32").Value = 2694.25
$ws.Range("K32").Value = 2694.25
$ws.Range("M32").Value = -2407.25

$ws.Range("H45").Value = 7995716
$ws.Range("I45").Value = 15985323
$ws.Range("K45").Value = 15985323
$ws.Range("M45").Value = -15984946

$ws.Range("H61").Value = 8502.105
$ws.Range("I61").Value = 8696.944
$ws.Range("K61").Value = 8696.944
$ws.Range("M61").Value = -8484.944

$ws.Range("H136").Value = 8502.105
$ws.Range("I136").Value = 8696.944
$ws.Range("K136").Value = 26090.832
$ws.Range("M136").Value = -23540.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5497405
$ws.Range("J99").Value = 3449.875
$ws.Range("L99").Value = 3449.875
$ws.Range("N99").Value = -6445.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 39179.6
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H61").Value = 39179.6
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H69").Value = 28047.223
$ws.Range("I69").Value = 19190.5
$ws.Range("K69").Value = 19190.5
$ws.Range("M69").Value = -18441.5

$ws.Range("H72").Value = 28047.223
$ws.Range("I72").Value = 19190.5
$ws.Range("K72").Value = 57571.5
$ws.Range("M72").Value = -53827.5

$ws.Range("H103").Value = 53413.43
$ws.Range("I103").Value = 53413.43
$ws.Range("K103").Value = 53413.43
$ws.Range("M103").Value = -52241.43

$ws.Range("H122").Value = 3796.6
$ws.Range("I122").Value = 3438.7144
$ws.Range("J122").Value = 4631.6665
$ws.Range("K122").Value = 10316.1432
$ws.Range("L122").Value = 13894.9995
$ws.Range("M122").Value = -7866.143199999999
$ws.Range("N122").Value = -18794.9995

$ws.Range("H132").Value = 69386.336
$ws.Range("I132").Value = 102647.2
$ws.Range("K132").Value = 307941.6
$ws.Range("M132").Value = -305411.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 898
$ws.Range("I8").Value = 898
$ws.Range("K8").Value = 2694
$ws.Range("M8").Value = -2555

$ws.Range("H15").Value = 913
$ws.Range("I15").Value = 479.66666
$ws.Range("J15").Value = 1173
$ws.Range("K15").Value = 1438.99998
$ws.Range("L15").Value = 3519
$ws.Range("M15").Value = -1298.99998
$ws.Range("N15").Value = -3799

$ws.Range("H23").Value = 218.63637
$ws.Range("J23").Value = 274.85715
$ws.Range("L23").Value = 824.5714499999999
$ws.Range("N23").Value = -1294.57145

$ws.Range("H68").Value = 1563.9656
$ws.Range("J68").Value = 1955
$ws.Range("L68").Value = 5865
$ws.Range("N68").Value = -7487

$ws.Range("H71").Value = 1563.9656
$ws.Range("J71").Value = 1955
$ws.Range("L71").Value = 17595
$ws.Range("N71").Value = -25707

$ws.Range("H107").Value = 834
$ws.Range("I107").Value = 740.4
$ws.Range("J107").Value = 951
$ws.Range("K107").Value = 2221.2
$ws.Range("L107").Value = 2853
$ws.Range("M107").Value = -301.1999999999998
$ws.Range("N107").Value = -6693

$ws.Range("H129").Value = 1177501.5
$ws.Range("I129").Value = 1334201.8
$ws.Range("J129").Value = 2249
$ws.Range("K129").Value = 4002605.4
$ws.Range("L129").Value = 6747
$ws.Range("M129").Value = -3997605.4
$ws.Range("N129").Value = -16747

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1758372.6
$ws.Range("I80").Value = 3499926.8
$ws.Range("K80").Value = 3499926.8
$ws.Range("M80").Value = -3498928.8

$ws.Range("H83").Value = 1758372.6
$ws.Range("I83").Value = 3499926.8
$ws.Range("K83").Value = 17499634
$ws.Range("M83").Value = -17494642

$ws.Range("H102").Value = 4533374
$ws.Range("I102").Value = 6537970
$ws.Range("K102").Value = 6537970
$ws.Range("M102").Value = -6536348

$ws.Range("H122").Value = 688545.4399999999
$ws.Range("I122").Value = 812008.25
$ws.Range("K122").Value = 2436024.75
$ws.Range("M122").Value = -2433574.75

$ws.Range("H135").Value = 70555.55499999999
$ws.Range("J135").Value = 70555.55499999999
$ws.Range("L135").Value = 70555.55499999999
$ws.Range("N135").Value = -80695.55499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H60").Value = 45000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H82").Value = 2316018.5
$ws.Range("I82").Value = 2779007.2
$ws.Range("K82").Value = 2779007.2
$ws.Range("M82").Value = -2778646.2

$ws.Range("H85").Value = 2316018.5
$ws.Range("I85").Value = 2779007.2
$ws.Range("K85").Value = 2779007.2
$ws.Range("M85").Value = -2777759.2

$ws.Range("H100").Value = 4054.077
$ws.Range("I100").Value = 2050.5
$ws.Range("J100").Value = 5771.4287
$ws.Range("K100").Value = 2050.5
$ws.Range("L100").Value = 5771.4287
$ws.Range("M100").Value = -1509.5
$ws.Range("N100").Value = -6853.4287

$ws.Range("H101").Value = 8021.25
$ws.Range("J101").Value = 8021.25
$ws.Range("L101").Value = 8021.25
$ws.Range("N101").Value = -14511.25

$ws.Range("H119").Value = 90499.5
$ws.Range("J119").Value = 90499.5
$ws.Range("L119").Value = 90499.5
$ws.Range("N119").Value = -100175.5

$ws.Range("H122").Value = 6973.909
$ws.Range("I122").Value = 4092.5
$ws.Range("K122").Value = 12277.5
$ws.Range("M122").Value = -9827.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10126.333
$ws.Range("J45").Value = 10126.333
$ws.Range("L45").Value = 10126.333
$ws.Range("N45").Value = -11108.333

$ws.Range("H107").Value = 50002644
$ws.Range("I107").Value = 66669828
$ws.Range("K107").Value = 200009484
$ws.Range("M107").Value = -200007564

$ws.Range("H136").Value = 5786.8
$ws.Range("I136").Value = 6874
$ws.Range("J136").Value = 2692.4614
$ws.Range("K136").Value = 20622
$ws.Range("L136").Value = 8077.3842
$ws.Range("M136").Value = -18072
$ws.Range("N136").Value = -13177.3842
